# Table Creator.xlsx edit
# - Rename table "Games" -> "Ts"
# - Remove the "Name" (varchar(50)) field row
# - Remove the "StartDate" (datetime2(3)) field row
# - Change the GetById stored-proc name formula (was "_Get")
# - Change the GetById parameter formula to use the PK column name directly
# - Widen column B, select C3 (cosmetic follow-on edits from the same session)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the table (C2 holds the table name used throughout the generated code)
$ws.Range("C2").Value2 = "Ts"

# Remove the two existing field rows (Name, StartDate) - ClearContents keeps
# the row's existing cell formatting/style intact, matching a Delete-key edit.
$ws.Range("B5:C5").ClearContents()
$ws.Range("B6:C6").ClearContents()

# Update the "Get" stored procedure name generator to "GetById"
$ws.Range("P2").Formula = '="CREATE PROCEDURE"&" ["&B2&"].[USP_"&C2&"_GetById]"'

# Update the GetById parameter declaration formula to reuse the PK column name
$ws.Range("P3").Formula = '=IF(NOT(ISBLANK(B4)),"@"&B4)'

# Widen column B (closest attainable width to the authored 16.86 chars) and
# leave the selection on C3, matching the end-of-session cursor position.
$ws.Columns.Item(2).ColumnWidth = 16

$ws.Range("C3").Select()
